# "adding averages and more checks"
#
# Training Dashboard: the PERIOD TO EXPIRE (H) / LAST UPDATE (I) columns were
# refreshed -- LAST UPDATE moved from 08-Sep-2025 to 16-Sep-2025 (8 days
# later), so every PERIOD TO EXPIRE value drops by 8.
#
# Exam Dashboard: the COMMENTS cell for the first exam got a more descriptive
# note, and its column was widened to fit.
#
# Header rows on both sheets get bold WHITE text (on the existing dark-blue
# fill) instead of plain bold black text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Training Dashboard"
$ws2 = $wb.Worksheets.Item(2)   # "Exam Dashboard"

# ---------------------------------------------------------------------------
# Training Dashboard: refresh LAST UPDATE date + recompute PERIOD TO EXPIRE
# ---------------------------------------------------------------------------
$lastUpdate = "'16-Sep-2025"   # leading apostrophe -> force literal text, not a date

$periodToExpire = @{
    3  = 351
    4  = 368
    5  = 336
    6  = 334
    7  = 357
    8  = 329
    9  = 355
    10 = 370
    11 = 332
    12 = 348
    13 = 338
    14 = 377
    15 = 79
    16 = -42
}

foreach ($row in 3..16) {
    $ws1.Range("H$row").Value = $periodToExpire[$row]
    $ws1.Range("I$row").Value = $lastUpdate
}

# ---------------------------------------------------------------------------
# Exam Dashboard: more descriptive comment + widen the COMMENTS column
# ---------------------------------------------------------------------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Columns(5).ColumnWidth = 14.17   # renders as column width 15

# ---------------------------------------------------------------------------
# Header styling: bold white text on both dashboards' header rows
# ---------------------------------------------------------------------------
$white = 16777215   # RGB(255,255,255)
$ws1.Range("A2:K2").Font.Color = $white
$ws2.Range("A2:G2").Font.Color = $white
